$wb = $excel.ActiveWorkbook

# Add the three new worksheets in order, after the last existing sheet
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$sMenuLinks = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$sMenuLinks.Name = "MenuLinks"
$sApplyNow = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sMenuLinks)
$sApplyNow.Name = "ApplyNow"
$sShopForYourCarNow = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sApplyNow)
$sShopForYourCarNow.Name = "ShopForYourCarNow"

# ---- MenuLinks sheet content ----
$menuLinks = @(
  "https://www.bankofamerica.com/auto-loans/",
  "https://www.bankofamerica.com/auto-loans/auto-refinance-loan/",
  "https://www.bankofamerica.com/auto-loans/auto-loan-rates/",
  "https://www.bankofamerica.com/auto-loans/disability-access-loans/",
  "https://www.bankofamerica.com/auto-loans/auto-loan-faq/",
  "https://www.bankofamerica.com/auto-loans/auto-loan-calculator/",
  "https://www.bankofamerica.com/auto-loans/auto-refinance-calculator/",
  "https://www.bankofamerica.com/auto-loans/how-car-loans-work/",
  "https://www.bankofamerica.com/auto-loans/financing-car/",
  "https://www.bankofamerica.com/auto-loans/when-to-refinance-a-car/",
  "https://www.bankofamerica.com/auto-loans/buying-new-or-used-cars/",
  "https://www.bankofamerica.com/auto-loans/lease-buyout/",
  "https://secure.bankofamerica.com/applynow/initialize-workflow.go?requesttype=VLSTATUS",
  "https://secure.bankofamerica.com/applynow/initialize-workflow.go?requesttype=SNR&flow=AUTO",
  "https://www.bankofamerica.com/customer-service/contact-us/auto-loans/?topicId=vehicle_oth_loans",
  "https://dealer-network.bankofamerica.com/"
)
for ($i = 0; $i -lt $menuLinks.Length; $i++) {
    $sMenuLinks.Cells.Item($i + 1, 1).Value = $menuLinks[$i]
}
$sMenuLinks.Columns.Item(1).ColumnWidth = 62.265625
$sMenuLinks.Range("A16").Select() | Out-Null

# ---- ApplyNow sheet content ----
$sApplyNow.Range("A1").Value = "Your Application"

# ---- ShopForYourCarNow sheet content ----
$sShopForYourCarNow.Range("A1").Value = "Shop for a Car and Financing in One Place at Bank of America"
$sShopForYourCarNow.Columns.Item(1).ColumnWidth = 54.9296875
$sShopForYourCarNow.Range("D4").Select() | Out-Null
$sShopForYourCarNow.Activate() | Out-Null

Write-Host "Done"
